$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 54 (shifts existing rows 54-186 down to 55-187,
# extending the used range to A1:R187), then populate the new row 54 with
# this week's new "Espinaca" price record.
$ws.Rows(54).Insert()

$ws.Range("A54").Value = 8
$ws.Range("B54").Value = "Terminal La Palmera de La Serena"
$ws.Range("C54").Value = "Coquimbo"
$ws.Range("D54").Value = 44526
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 100112012
$ws.Range("G54").Value = "Espinaca"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 3320
$ws.Range("K54").Value = 400
$ws.Range("L54").Value = 500
$ws.Range("M54").Value = 450
$ws.Range("N54").Value = "$/atado 300 a 500 gramos"
$ws.Range("O54").Value = "Provincia del Elquí"
$ws.Range("P54").Value = 900
$ws.Range("Q54").Value = 0.5
$ws.Range("R54").Value = "Hortaliza"
